# Moved statistical datasets and results
# - Inserts a new column at A (pushing the original A:D data to B:E)
# - Inserts a new row at 1 (pushing the data down, for a header row)
# - Adds header labels in B1:E1 ("Valid","T","Z","p-value")
# - Adds row labels in A2:A22 describing each Wilcoxon test pairing
# - Restores / sets the column widths to match the new layout

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column A, and a new row before row 1.
$ws.Columns.Item(1).Insert()
$ws.Rows.Item(1).Insert()

# 2. Header row (row 1) - columns B:E
$ws.Range("B1").Value = "Valid"
$ws.Range("C1").Value = "T"
$ws.Range("D1").Value = "Z"
$ws.Range("E1").Value = "p-value"

# 3. Row labels (column A, rows 2:22)
$labels = @(
    "CyclomaticComplexity(CC) & CyclomaticComplexity(CC)",
    "MaintainabilityIndex & MaintainabilityIndex",
    "NbUniqueOperands & NbUniqueOperands",
    "NbUniqueOperands & EffortToImplement",
    "NbOperands & NbOperands",
    "NbOperands & EffortToImplement",
    "NbUniqueOperators & NbUniqueOperators",
    "NbUniqueOperators & EffortToImplement",
    "NbOperators & NbOperators",
    "ProgramLength & ProgramLength",
    "ProgramLength & EffortToImplement",
    "VocabularySize & VocabularySize",
    "ProgramVolume & ProgramVolume",
    "DifficultyLevel & DifficultyLevel",
    "ProgramLevel & ProgramLevel",
    "EffortToImplement & NbUniqueOperands",
    "EffortToImplement & NbOperands",
    "EffortToImplement & NbUniqueOperators",
    "EffortToImplement & ProgramLength",
    "EffortToImplement & EffortToImplement",
    "TimeToImplement & TimeToImplement"
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $labels[$i]
}

# 4. Column width for the newly added column A (columns B:E keep the
#    widths they already had before the insert, carried over automatically).
$ws.Columns.Item(1).ColumnWidth = 53.7
